$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for each data row (2..11), columns D,L,M,N,O,P,Q,R,S,T
# (values correspond to re-sorting the original rows by date ascending)

$rows = @(
    @{ Row = 2;  D = 44272; L = "Primera"; M = 100; N = 9000;  O = 10000; P = 9500;  Q = "`$/caja 15 kilos granel";    R = "Región de O'Higgins"; S = 633; T = 15 }
    @{ Row = 3;  D = 44272; L = "Segunda"; M = 50;  N = 8000;  O = 8000;  P = 8000;  Q = "`$/caja 15 kilos granel";    R = "Región de O'Higgins"; S = 533; T = 15 }
    @{ Row = 4;  D = 44299; L = "Primera"; M = 100; N = 10000; O = 11000; P = 10500; Q = "`$/caja 18 kilos granel";    R = "Región del Maule";    S = 583; T = 18 }
    @{ Row = 5;  D = 44299; L = "Segunda"; M = 50;  N = 9000;  O = 9000;  P = 9000;  Q = "`$/caja 18 kilos granel";    R = "Región del Maule";    S = 500; T = 18 }
    @{ Row = 6;  D = 44307; L = "Primera"; M = 50;  N = 10000; O = 10000; P = 10000; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 556; T = 18 }
    @{ Row = 7;  D = 44307; L = "Segunda"; M = 50;  N = 8000;  O = 8000;  P = 8000;  Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 444; T = 18 }
    @{ Row = 8;  D = 44316; L = "Primera"; M = 100; N = 9000;  O = 10000; P = 9500;  Q = "`$/caja 18 kilos granel";    R = "Región de O'Higgins"; S = 528; T = 18 }
    @{ Row = 9;  D = 44363; L = "Primera"; M = 100; N = 9000;  O = 10000; P = 9500;  Q = "`$/caja 15 kilos empedrada"; R = "Región de O'Higgins"; S = 633; T = 15 }
    @{ Row = 10; D = 44358; L = "Primera"; M = 100; N = 11000; O = 12000; P = 11500; Q = "`$/caja 18 kilos granel";    R = "Región de O'Higgins"; S = 639; T = 18 }
    @{ Row = 11; D = 44425; L = "Primera"; M = 100; N = 12000; O = 13000; P = 12500; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 694; T = 18 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($rowNum, 12).Value = $r.L   # L: Calidad
    $ws.Cells.Item($rowNum, 13).Value = $r.M   # M: Volumen
    $ws.Cells.Item($rowNum, 14).Value = $r.N   # N: Precio mínimo
    $ws.Cells.Item($rowNum, 15).Value = $r.O   # O: Precio máximo
    $ws.Cells.Item($rowNum, 16).Value = $r.P   # P: Precio promedio ponderado
    $ws.Cells.Item($rowNum, 17).Value = $r.Q   # Q: Unidad de comercialización
    $ws.Cells.Item($rowNum, 18).Value = $r.R   # R: Origen
    $ws.Cells.Item($rowNum, 19).Value = $r.S   # S: Precio $/Kg
    $ws.Cells.Item($rowNum, 20).Value = $r.T   # T: Kg / unidad
}
